# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" note (A1) with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newline = [char]10
$text = "Conversión del día 💰" + $newline +
        "✅ Dólar paralelo: 68" + $newline +
        $newline +
        "Binance" + $newline +
        "✅ 1000 Bs = 12.5 = 49866.24 pesos" + $newline +
        "✅ 49866.24 pesos = 12.42 = 960.49 Bs" + $newline +
        $newline +
        "Promedio competencia" + $newline +
        "✅ Tasa pesos: 20" + $newline +
        "✅ Tasa Bs: 20" + $newline +
        "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $text

# --- Update the "tasas" sheet numeric rates ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 79.994
$wsTasas.Range("O10").Value = 3989
$wsTasas.Range("N12").Value = 4014
$wsTasas.Range("O12").Value = 77.315
